# Notified_Production_Wind.xlsx update:
#  - Column A (timestamps, rows 2-97): shift every date serial forward by
#    exactly 7 days (one week later - "switching to winter DST time").
#  - Column B (values, rows 2-93): replace with the newly fetched readings.
#    Rows 94-97 keep their existing value of 0 (only the date in column A
#    changes for those rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column B values for rows 2 through 93 (92 values, in row order).
$newValues = @(
    1199.288, 1199.726, 1196.317, 1195.015, 1202.108, 1212.43, 1221.562, 1228.64, 1267.47, 1279.135,
    1290.313, 1300.841, 1376.972, 1388.671, 1406.527, 1422.275, 1505.394, 1526.395, 1539.993, 1548.989,
    1623.143, 1635.552, 1649.058, 1657.691, 1700.383, 1700.858, 1699.654, 1700.248, 1787.311, 1801.938,
    1800.569, 1826.572, 1876.659, 1885.045, 1874.511, 1884.377, 1939.765, 1939.032, 1936.875, 1933.973,
    1903.804, 1896.775, 1886.578, 1879.451, 1863.773, 1862.55, 1867.589, 1858.665, 1837.728, 1835.666,
    1830.982, 1822.491, 1815.844, 1828.701, 1831.139, 1830.163, 1839.614, 1845.722, 1853.768, 1862.423,
    1950.152, 1961.682, 1981.396, 1997.301, 2128.534, 2142.771, 2152.226, 2155.577, 2205.565, 2212.186,
    2216.976, 2219.762, 2188.144, 2181.107, 2176.217, 2170.559, 2087.455, 2078.662, 2070.376, 2062.612,
    1971.327, 1954.78, 1914.862, 1896.872, 1758.969, 1735.241, 1709.694, 1683.212, 1537.174, 1515.99,
    1494.493, 1464.91
)

# Shift every timestamp in column A (rows 2-97) forward by 7 days.
for ($row = 2; $row -le 97; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value2 = $cellA.Value2 + 7
}

# Overwrite column B (rows 2-93) with the new readings.
for ($row = 2; $row -le 93; $row++) {
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row - 2]
}
